$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

$ws.Range("H7").Value = "Number of colonies * dilution factor"
$ws.Range("B7").Value = "Initial colony count"
$ws.Range("N7").Value = "Number per 100μl * aliquot volume/well volume"
$n7 = $ws.Range("N7")
$chars = $n7.Characters(15, 31)
$chars.Font.Size = 11
$chars.Font.Name = "Calibri"
$chars.Font.ColorIndex = 1

$ws.Range("O3:Q6").Formula = "=I3*10"

$ws.Range("N8").Select()
